# Add data loader module:
#  - Translate header row to Chinese, repurpose column D as a "link" column,
#    and move the job description text to a new column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# E1 is a brand-new header cell; give it the same header style (bold,
# centered, bordered) as the existing header cells before setting its text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A1").Value = "公司名字"
$ws.Range("B1").Value = "岗位名"
$ws.Range("C1").Value = "地点"
$ws.Range("D1").Value = "link"
$ws.Range("E1").Value = "job description"

# --- Move existing job descriptions from D to E, then fill D with links ---
$descriptions = @{
    2 = "We are seeking a Software Engineer Intern to join our backend development team. The intern will work on designing and implementing scalable web services using Java and Spring Boot. Responsibilities include developing RESTful APIs, working with databases, and collaborating with senior engineers on system architecture. No citizenship requirements. This is an entry-level internship position perfect for students."
    3 = "Looking for a Machine Learning Engineer Intern to help build and deploy AI models in production. The role involves working with Python, TensorFlow, and AWS to create recommendation systems and computer vision applications. The intern will collaborate with data scientists and software engineers to optimize model performance. This internship is open to all students regardless of citizenship status. No prior industry experience required."
    4 = "Quantitative Research Intern position available for students interested in algorithmic trading and financial modeling. The role requires strong mathematical background and programming skills in Python or R. Intern will develop trading strategies, perform statistical analysis, and work with large financial datasets. This position requires US citizenship or permanent residency due to regulatory compliance requirements."
}

$links = @{
    2 = "https://techcorp.com/careers/swe-intern"
    3 = "https://aisolutions.com/jobs/mle-intern"
    4 = "https://dataflow.com/careers/quant-intern"
}

foreach ($row in 2..4) {
    $ws.Cells.Item($row, 5).Value = $descriptions[$row]
    $ws.Cells.Item($row, 4).Value = $links[$row]
}
